$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Sheet1: re-point the 3 data rows to the new row order / new price values
#    Row order becomes: Helios Amante..., Midori Bamboo..., Addison...
# ---------------------------------------------------------------------------

# Column A holds plain text (never looks numeric), so .Value is safe.
$ws1.Range("A2").Value = "Helios Amante Mango Wood Book Shelf - Brown"
$ws1.Range("A3").Value = "Midori Bamboo 5-Tier Book Shelf - Light Brown"
$ws1.Range("A4").Value = "Addison 3-Tier Book Shelf - Beige"

# Column B holds numeric-looking text ("8,999" etc.) that must stay TEXT
# (shared string), matching the source workbook. Assigning such a string
# straight to .Value makes Excel "smart" parse it into a real number, so we
# stage the literal text in a scratch range that is explicitly formatted as
# Text ("@"), copy it, and paste-special *values only* into the target
# cells (this carries the text-ness across without carrying the "@" format
# onto the destination cells). The scratch column is then deleted outright
# so no trace of it remains in the saved sheet.
$scratch = $ws1.Range("Z1:Z3")
$scratch.NumberFormat = "@"
$ws1.Range("Z1").Value = "8,999"
$ws1.Range("Z2").Value = "4,997"
$ws1.Range("Z3").Value = "4,499"

$ws1.Range("Z1").Copy()
$ws1.Range("B2").PasteSpecial(-4163)
$ws1.Range("Z2").Copy()
$ws1.Range("B3").PasteSpecial(-4163)
$ws1.Range("Z3").Copy()
$ws1.Range("B4").PasteSpecial(-4163)

$ws1.Columns.Item(26).Delete()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Sheet2 ("Home Accessories Items" category list), inserted right after
#    Sheet1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

$ws2.Range("A1").Value = "Home Accessories Items"
$ws2.Range("A2").Value = "Candle Holders And Lanterns"
$ws2.Range("A3").Value = "Clocks"
$ws2.Range("A4").Value = "Desk Accessories"
$ws2.Range("A5").Value = "Figurines"
$ws2.Range("A6").Value = "Gift Accessories"
$ws2.Range("A7").Value = "Lamps"
$ws2.Range("A8").Value = "Mirrors"
$ws2.Range("A9").Value = "Ornaments"
$ws2.Range("A10").Value = "Potpouri Bowl and Platter"

# Re-use Sheet1's header style (bold font + themed fill + centered) for the
# title cell instead of building a brand-new style.
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Columns.Item(1).ColumnWidth = 40.8

$ws2.Activate()
$ws2.Range("G9").Select()

# ---------------------------------------------------------------------------
# 3. Sheet3 ("Login Error Message" copy), inserted right after Sheet2.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)

$ws3.Range("A1").Value = "Login Error Message"
$ws3.Range("A2").Value = "Email must be a valid email address."

$ws1.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Columns.Item(1).ColumnWidth = 42.65

$ws3.Activate()
$ws3.Range("C3").Select()

# ---------------------------------------------------------------------------
# 4. Restore Sheet1 as the active / selected sheet.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()

Write-Host "Sheet count:" $wb.Worksheets.Count
Write-Host "A2:" $ws1.Range("A2").Value2 " B2:" $ws1.Range("B2").Value2
Write-Host "A3:" $ws1.Range("A3").Value2 " B3:" $ws1.Range("B3").Value2
Write-Host "A4:" $ws1.Range("A4").Value2 " B4:" $ws1.Range("B4").Value2
